$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($ws, $addr, $val) {
    # Force the target cell to Text format before assigning so numeric-
    # looking strings (e.g. "3.28", "1.00") are kept as literal text
    # instead of being auto-coerced to a number by the Value setter.
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    # Reset the style back to Normal/default afterwards so we don't leave
    # a lingering text-format style on cells that didn't have one before.
    $cell.Style = "Normal"
}


Set-CellText $ws 'D2' '30.479.37'
Set-CellText $ws 'E2' '  +1.80%  '
Set-CellText $ws 'D3' '1.673.04'
Set-CellText $ws 'E3' '  +2.39%  '
Set-CellText $ws 'E4' '  -0.06%  '
Set-CellText $ws 'D5' '219.49'
Set-CellText $ws 'E5' '  +2.33%  '
Set-CellText $ws 'E6' '  +2.42%  '
Set-CellText $ws 'D7' '1.00'
Set-CellText $ws 'E7' '  +0.02%  '
Set-CellText $ws 'D8' '29.64'
Set-CellText $ws 'E8' '  +3.85%  '
Set-CellText $ws 'E9' '  +2.55%  '
Set-CellText $ws 'D10' '0.0638'
Set-CellText $ws 'E10' '  +5.00%  '
Set-CellText $ws 'E11' '  -0.52%  '
Set-CellText $ws 'D12' '1.913.26'
Set-CellText $ws 'E12' '  +2.43%  '
Set-CellText $ws 'D13' '0.614'
Set-CellText $ws 'E13' '  +9.08%  '
Set-CellText $ws 'D14' '1.664.51'
Set-CellText $ws 'E14' '  +1.91%  '
Set-CellText $ws 'D15' '10.19'
Set-CellText $ws 'E15' '  +8.73%  '
Set-CellText $ws 'D16' '3.98'
Set-CellText $ws 'E16' '  +3.27%  '
Set-CellText $ws 'D17' '30.521.98'
Set-CellText $ws 'E17' '  +1.91%  '
Set-CellText $ws 'E18' '  +3.60%  '
Set-CellText $ws 'D19' '242.55'
Set-CellText $ws 'E19' '  +0.04%  '
Set-CellText $ws 'D20' '0.0₃0720'
Set-CellText $ws 'E20' '  +2.74%  '
Set-CellText $ws 'E21' '  -0.06%  '
Set-CellText $ws 'E22' '  +3.05%  '
Set-CellText $ws 'D23' '9.97'
Set-CellText $ws 'E23' '  +0.90%  '
Set-CellText $ws 'E24' '  +0.09%  '
Set-CellText $ws 'D25' '158.37'
Set-CellText $ws 'E25' '  +0.54%  '
Set-CellText $ws 'D26' '15.85'
Set-CellText $ws 'E26' '  +2.09%  '
Set-CellText $ws 'D27' '0.112'
Set-CellText $ws 'E27' '  +2.32%  '
Set-CellText $ws 'D28' '6.67'
Set-CellText $ws 'E28' '  +0.84%  '
Set-CellText $ws 'D29' '1.00'
Set-CellText $ws 'E29' '  -0.02%  '
Set-CellText $ws 'D30' '0.0495'
Set-CellText $ws 'E30' '  +1.81%  '
Set-CellText $ws 'D31' '1.14'
Set-CellText $ws 'E31' '  +2.64%  '
Set-CellText $ws 'E32' '  +2.57%  '
Set-CellText $ws 'B33' 'Maker'
Set-CellText $ws 'C33' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-CellText $ws 'D33' '1.500.31'
Set-CellText $ws 'E33' '  +5.26%  '
Set-CellText $ws 'B34' 'InternetComputer(DFINITY)'
Set-CellText $ws 'C34' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-CellText $ws 'D34' '3.28'
Set-CellText $ws 'E34' '  +3.17%  '
Set-CellText $ws 'E35' '  +6.93%  '
Set-CellText $ws 'D36' '84.29'
Set-CellText $ws 'E36' '  +10.69%  '
Set-CellText $ws 'E37' '  -0.85%  '
Set-CellText $ws 'D38' '0.598'
Set-CellText $ws 'E38' '  +8.31%  '
Set-CellText $ws 'D39' '0.0177'
Set-CellText $ws 'E39' '  +4.94%  '
Set-CellText $ws 'E40' '  -4.58%  '
Set-CellText $ws 'E41' '  -0.10%  '
Set-CellText $ws 'D42' '0.838'
Set-CellText $ws 'E42' '  +1.35%  '
Set-CellText $ws 'D43' '1.97'
Set-CellText $ws 'E43' '  -1.33%  '
Set-CellText $ws 'E44' '  +1.42%  '
Set-CellText $ws 'E45' '  +0.47%  '
Set-CellText $ws 'D46' '1.00'
Set-CellText $ws 'E46' '  +0.01%  '
Set-CellText $ws 'D47' '5.55'
Set-CellText $ws 'E47' '  +3.33%  '
Set-CellText $ws 'D48' '50.91'
Set-CellText $ws 'E48' '  -4.03%  '
Set-CellText $ws 'D49' '1.804.68'
Set-CellText $ws 'E49' '  +1.64%  '
Set-CellText $ws 'D50' '94.61'
Set-CellText $ws 'E50' '  +4.60%  '
Set-CellText $ws 'E51' '  +0.67%  '
